# Atualização automática de preços de eletricidade
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 46069

$ws.Range("B2").Value = 0.01
$ws.Range("C2").Value = 0
$ws.Range("D2").Value = -0.01
$ws.Range("E2").Value = -0.02
$ws.Range("F2").Value = -0.03
$ws.Range("G2").Value = -0.01
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0.03
$ws.Range("J2").Value = 0.1
$ws.Range("K2").Value = 0.34
$ws.Range("L2").Value = 0.1
$ws.Range("M2").Value = 0
$ws.Range("N2").Value = 0
$ws.Range("O2").Value = -0.01
$ws.Range("P2").Value = 0
$ws.Range("Q2").Value = 0
$ws.Range("R2").Value = 0
$ws.Range("S2").Value = 1.48
$ws.Range("T2").Value = 6.69
$ws.Range("U2").Value = 24.35
$ws.Range("V2").Value = 24.75
$ws.Range("W2").Value = 13.84
$ws.Range("X2").Value = 10.16
$ws.Range("Y2").Value = 2.18
$ws.Range("Z2").Value = 3.5

$ws.Range("AB2").Value = 12.73
$ws.Range("AD2").Value = 19.3
$ws.Range("AE2").Value = "18h-20h"
$ws.Range("AF2").Value = 15.52
$ws.Range("AG2").Value = "0h-23h"
